$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.819.48"
$ws.Range("E2").Value = "  -1.98%  "
$ws.Range("D3").Value = "1.801.08"
$ws.Range("E3").Value = "  -1.35%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "308.79"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.70%  "
$ws.Range("E6").Value = "  +0.07%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4659"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +4.30%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3678"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07348"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.83%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8667"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.16%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.36"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.25%  "
$ws.Range("D12").Value = "1.868.42"
$ws.Range("E12").Value = "  +2.28%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.347"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.39%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.510"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.64%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.07025"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.24%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "91.28"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.71%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.001"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.02%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008697"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.89%  "
$ws.Range("E19").Value = "  +0.13%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.63"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.49%  "
$ws.Range("D21").Value = "26.823.65"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.289"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.32%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.60"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.10%  "
$ws.Range("D24").Value = "2.078.45"
$ws.Range("E24").Value = "  +0.90%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.903"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.26%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "151.29"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.08%  "
$ws.Range("E27").Value = "  -1.82%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.116"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -8.27%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.234"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.71%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "115.81"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.74%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08908"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.39%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7542"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.15%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.930"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.56%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.145"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.17%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.441"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.70%  "
$ws.Range("E36").Value = "  +0.08%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.103"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.85%  "
$ws.Range("E38").Value = "  -2.10%  "
$ws.Range("E39").Value = "  -1.15%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.932"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.41%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.182"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.85%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5264"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.65%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.330"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.44%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1658"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.52%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.450"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.53%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4994"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.47%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.26"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.29%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "103.96"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.29%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.000"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.06%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.663"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.01%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06281"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.89%  "
